$wb = $excel.ActiveWorkbook

# The int-constant table: rename "MaxTotalSkillLevel" -> "MaxTotalSpellLevel"
# and change its value from 100 to a temporary 10.
$ws1 = $wb.Worksheets.Item("GlobalConstantIntTable")
$ws1.Range("A21").Value = "MaxTotalSpellLevel"
$ws1.Range("B21").Value = 10

# Make the int-constant sheet the active sheet/tab, with A2 selected.
$ws1.Activate()
$ws1.Range("A2").Select()
